{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate \"This is the first commit I have made to this project.\" \u2014 the\n// paragraph the new \"modifies\" paragraph must be inserted directly after.\nconst anchorText = \"This is the first commit I have made to this project.\";\nlet anchorParagraph = paragraphs.items.find((p) => p.text === anchorText);\nif (!anchorParagraph) {\n  // Fall back to the first paragraph if the exact text can't be matched.\n  anchorParagraph = paragraphs.items[0];\n}\n\n// Insert a new paragraph with text \"modifies\" right after it, and before the\n// existing (originally second) empty paragraph.\nanchorParagraph.insertParagraph(\"modifies\", \"After\");\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate \"This is the first commit I have made to this project.\" \u2014 the\n# paragraph the new \"modifies\" paragraph must be inserted directly after.\n$anchorText = \"This is the first commit I have made to this project.\"\n$anchorParagraph = $null\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text.TrimEnd([char]13, [char]7)\n    if ($t -eq $anchorText) {\n        $anchorParagraph = $p\n        break\n    }\n}\nif ($anchorParagraph -eq $null) {\n    # Fall back to the first paragraph if the exact text can't be matched.\n    $anchorParagraph = $d.Paragraphs.Item(1)\n}\n\n# Insert a brand-new (empty) paragraph right after it, and before the\n# existing (originally second) empty paragraph.\n$anchorParagraph.Range.InsertParagraphAfter() | Out-Null\n\n# The freshly inserted paragraph directly follows the anchor - fill it with \"modifies\".\n$newIndex = $anchorParagraph.Index + 1\n$newParagraph = $d.Paragraphs.Item($newIndex)\n$newParagraph.Range.InsertAfter(\"modifies\")\n"}
